$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data cells on row 3: a date+time, a percentage, and a date.
# Custom number formats get registered (and their cellXfs entries appended)
# in the order the cells below are first formatted, so this order is chosen
# to line up the resulting style/format table with the target workbook.

# C3: date+time value (serial 42309.604166666664 == 2015-11-01 14:30), custom
# date/time format.
$ws.Range("C3").Value = 42309.604166666664
$ws.Range("C3").NumberFormat = "m/d/yy\ h:mm;@"

# D3: percentage value, custom percentage format with 4 decimal places.
$ws.Range("D3").Value = 0.91323449999999995
$ws.Range("D3").NumberFormat = "0.0000%"

# B3: date-only value (serial 42309 == 2015-11-01), custom date format.
$ws.Range("B3").Value = 42309
$ws.Range("B3").NumberFormat = "[$-409]mmmm\ d\,\ yyyy;@"

# Move the cursor/selection onto the newly-populated cell (matches the
# author's recorded cursor position after making the edit).
$null = $ws.Range("B3").Select()
